# Insert a new weekly record for "Pepino ensalada" (Macroferia Regional de Talca)
# right before the existing row 408. This pushes the old rows 408-439 down to
# 409-440 (dimension grows from R439 to R440) and populates the newly opened
# row 408 with the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 408:439 down by one row, opening up a blank row 408.
$ws.Rows.Item(408).Insert()

# Fill in the new row 408 with the new weekly data point.
$ws.Range("A408").Value = 5
$ws.Range("B408").Value = "Macroferia Regional de Talca"
$ws.Range("C408").Value = "Maule"
$ws.Range("D408").Value = 44746
$ws.Range("E408").Value = 7
$ws.Range("F408").Value = 100112043
$ws.Range("G408").Value = "Pepino ensalada"
$ws.Range("H408").Value = "Sin especificar"
$ws.Range("I408").Value = "Primera"
$ws.Range("J408").Value = 400
$ws.Range("K408").Value = 18000
$ws.Range("L408").Value = 18000
$ws.Range("M408").Value = 18000
$ws.Range("N408").Value = "`$/caja 60 unidades"
$ws.Range("O408").Value = "Región de Arica y Parinacota"
$ws.Range("P408").Value = 300
$ws.Range("Q408").Value = 60
$ws.Range("R408").Value = "Hortaliza"
